# Update cryptocurrency price/volume figures (cryptos.xlsx refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ('Price') cells hold plain numeric-looking text (e.g. '1.00', '563.32').
# A leading apostrophe forces Excel to keep them as text instead of coercing to
# a Double (which would drop the trailing zero / add floating point noise).

$ws.Range('D2').Value = '66.194.67'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.315.21'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''563.32'
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('D6').Value = '''186.04'
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '3.308.57'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('D9').Value = '''0.575'
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('E10').Value = '  -4.97%  '
$ws.Range('D11').Value = '''0.575'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('D12').Value = '''46.09'
$ws.Range('E12').Value = '  -2.79%  '
$ws.Range('D13').Value = '''0.0000265'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('D14').Value = '3.844.89'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('D16').Value = '''590.86'
$ws.Range('E16').Value = '  -8.37%  '
$ws.Range('D17').Value = '66.095.16'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D19').Value = '3.312.91'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').Value = '''17.68'
$ws.Range('E20').Value = '  -2.20%  '
$ws.Range('E21').Value = '  -4.03%  '
$ws.Range('D22').Value = '''0.896'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('D23').Value = '''18.11'
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').Value = '''5.03'
$ws.Range('E24').Value = '  +2.20%  '
$ws.Range('D25').Value = '''98.40'
$ws.Range('E25').Value = '  -9.04%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').Value = '''9.41'
$ws.Range('E28').Value = '  -1.67%  '
$ws.Range('D29').Value = '''8.46'
$ws.Range('E29').Value = '  -2.69%  '
$ws.Range('E30').Value = '  +0.72%  '
$ws.Range('D31').Value = '''6.67'
$ws.Range('E31').Value = '  +5.79%  '
$ws.Range('E32').Value = '  -6.25%  '
$ws.Range('D33').Value = '''562.87'
$ws.Range('E33').Value = '  +8.27%  '
$ws.Range('D34').Value = '''10.84'
$ws.Range('E34').Value = '  -1.98%  '
$ws.Range('D35').Value = '3.791.64'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  -1.60%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('D38').Value = '''56.02'
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').Value = '''33.35'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('D41').Value = '0.0₃0687'
$ws.Range('E41').Value = '  -6.78%  '
$ws.Range('E42').Value = '  -7.54%  '
$ws.Range('E43').Value = '  +3.86%  '
$ws.Range('E44').Value = '  -5.17%  '
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '''3.07'
$ws.Range('E47').Value = '  -8.84%  '
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('D49').Value = '''1.00'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('E50').Value = '  -3.01%  '
$ws.Range('D51').Value = '''128.06'
$ws.Range('E51').Value = '  +4.98%  '
